# Apply strikethrough formatting to the paragraphs describing items that
# were dropped from the requirements ("Throwing in the towel").
#
# The following paragraphs (identified by their exact trimmed text) get
# struck through, matching Word's behavior when the paragraph is selected
# and the Strikethrough toolbar button is toggled on:
#   - Primary Contact: Create automation to always have the important
#     contacts on the job application record.
#   - Clean Up Stale Jobs Applications
#   - Create an asynchronous process that checks if a job application is
#     stale and moves the record status to closed. Update the notes field
#     that the job application was closed by an automated process.
#   - Stale Criteria:
#   - Status is not Closed or Accepted
#   - Follow-up Date 30 days old or more

$d = $word.ActiveDocument

$targets = @(
    "Primary Contact: Create automation to always have the important contacts on the job application record.",
    "Clean Up Stale Jobs Applications",
    "Create an asynchronous process that checks if a job application is stale and moves the record status to closed. Update the notes field that the job application was closed by an automated process.",
    "Stale Criteria:",
    "Status is not Closed or Accepted",
    "Follow-up Date 30 days old or more"
)

foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text.TrimEnd([char]13, [char]7)
    foreach ($target in $targets) {
        if ($t -eq $target) {
            $p.Range.Font.StrikeThrough = 1
        }
    }
}
